$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("AJ9").Value = 9.5
$ws.Range("Q9").Value = 1.83
$ws.Range("R9").Value = 1.98

# Row 12
$ws.Range("AA12").Value = 9
$ws.Range("AB12").Value = 8.5
$ws.Range("AD12").Value = 23
$ws.Range("AF12").Value = 12
$ws.Range("AG12").Value = 21
$ws.Range("AI12").Value = 251
$ws.Range("G12").Value = 1.27
$ws.Range("H12").Value = 6
$ws.Range("I12").Value = 9.5
$ws.Range("J12").Value = 1.67
$ws.Range("L12").Value = 8
$ws.Range("N12").Value = 21
$ws.Range("W12").Value = 1.91
$ws.Range("X12").Value = 1.91
$ws.Range("Y12").Value = 9.5

# Row 13
$ws.Range("AC13").Value = 23
$ws.Range("AD13").Value = 29
$ws.Range("AE13").Value = 12
$ws.Range("AF13").Value = 7
$ws.Range("AG13").Value = 13
$ws.Range("AH13").Value = 41
$ws.Range("AJ13").Value = 8.5
$ws.Range("H13").Value = 3.6
$ws.Range("Q13").Value = 1.88
$ws.Range("R13").Value = 1.98
$ws.Range("U13").Value = 1.36
$ws.Range("V13").Value = 3
$ws.Range("W13").Value = 1.7
$ws.Range("X13").Value = 2.05

# Row 14
$ws.Range("AB14").Value = 51
$ws.Range("AE14").Value = 5.5
$ws.Range("AF14").Value = 6.5
$ws.Range("AJ14").Value = 4.75
$ws.Range("AK14").Value = 7
$ws.Range("AL14").Value = 10
$ws.Range("AM14").Value = 15
$ws.Range("AN14").Value = 21
$ws.Range("AP14").Value = 2.05
$ws.Range("AQ14").Value = 1.75
$ws.Range("G14").Value = 5
$ws.Range("H14").Value = 3
$ws.Range("I14").Value = 1.8
$ws.Range("J14").Value = 6
$ws.Range("K14").Value = 1.91
$ws.Range("L14").Value = 2.6
$ws.Range("M14").Value = 1.13
$ws.Range("N14").Value = 6
$ws.Range("O14").Value = 1.57
$ws.Range("P14").Value = 2.25
$ws.Range("Q14").Value = 2.88
$ws.Range("R14").Value = 1.4
$ws.Range("S14").Value = 6
$ws.Range("T14").Value = 1.13
$ws.Range("U14").Value = 1.62
$ws.Range("V14").Value = 2.2
$ws.Range("W14").Value = 2.5
$ws.Range("X14").Value = 1.5
$ws.Range("Y14").Value = 9.5
$ws.Range("Z14").Value = 23

# Row 16
$ws.Range("AC16").Value = 21
$ws.Range("AE16").Value = 6
$ws.Range("AF16").Value = 6
$ws.Range("AG16").Value = 19
$ws.Range("AH16").Value = 67
$ws.Range("AJ16").Value = 8.5
$ws.Range("AP16").Value = 1.98
$ws.Range("AQ16").Value = 1.83
$ws.Range("G16").Value = 2.1
$ws.Range("H16").Value = 2.8
$ws.Range("I16").Value = 3.6
$ws.Range("J16").Value = 3
$ws.Range("K16").Value = 1.91
$ws.Range("M16").Value = 1.13
$ws.Range("N16").Value = 6
$ws.Range("O16").Value = 1.5
$ws.Range("P16").Value = 2.5
$ws.Range("Q16").Value = 2.6
$ws.Range("R16").Value = 1.48
$ws.Range("S16").Value = 5.5
$ws.Range("T16").Value = 1.14
$ws.Range("U16").Value = 1.57
$ws.Range("V16").Value = 2.25
$ws.Range("W16").Value = 2.2
$ws.Range("X16").Value = 1.62

# Row 18
$ws.Range("AB18").Value = 10
$ws.Range("AF18").Value = 8.5
$ws.Range("AG18").Value = 21
$ws.Range("AH18").Value = 67
$ws.Range("AJ18").Value = 15
$ws.Range("AK18").Value = 29
$ws.Range("AL18").Value = 19
$ws.Range("AM18").Value = 67
$ws.Range("AN18").Value = 51
$ws.Range("AO18").Value = 51
$ws.Range("AP18").Value = 1.43
$ws.Range("AQ18").Value = 2.85
$ws.Range("AR18").Value = 2.5
$ws.Range("AS18").Value = 1.53
$ws.Range("G18").Value = 1.48
$ws.Range("H18").Value = 4.1
$ws.Range("I18").Value = 5.75
$ws.Range("J18").Value = 2.05
$ws.Range("K18").Value = 2.25
$ws.Range("L18").Value = 6.5
$ws.Range("M18").Value = 1.04
$ws.Range("N18").Value = 13
$ws.Range("O18").Value = 1.25
$ws.Range("P18").Value = 3.75
$ws.Range("Q18").Value = 1.88
$ws.Range("R18").Value = 1.98
$ws.Range("S18").Value = 3.25
$ws.Range("T18").Value = 1.33
$ws.Range("W18").Value = 2
$ws.Range("X18").Value = 1.73
$ws.Range("Z18").Value = 6.5

# Row 23
$ws.Range("AB23").Value = 14
$ws.Range("AC23").Value = 14
$ws.Range("AD23").Value = 26
$ws.Range("AE23").Value = 7.5
$ws.Range("AF23").Value = 7
$ws.Range("AG23").Value = 15
$ws.Range("AH23").Value = 70
$ws.Range("AI23").Value = 500
$ws.Range("AJ23").Value = 12
$ws.Range("AK23").Value = 23
$ws.Range("G23").Value = 1.75
$ws.Range("H23").Value = 3.55
$ws.Range("I23").Value = 4.05
$ws.Range("J23").Value = 2.35
$ws.Range("K23").Value = 2.18
$ws.Range("M23").Value = 1.06
$ws.Range("N23").Value = 7.5
$ws.Range("O23").Value = 1.28
$ws.Range("P23").Value = 3.35
$ws.Range("Q23").Value = 1.83
$ws.Range("R23").Value = 1.87
$ws.Range("S23").Value = 3
$ws.Range("T23").Value = 1.34
$ws.Range("U23").Value = 1.39
$ws.Range("V23").Value = 2.77
$ws.Range("W23").Value = 1.78
$ws.Range("X23").Value = 1.93
$ws.Range("Y23").Value = 7.2
$ws.Range("Z23").Value = 8.25

# Row 27
$ws.Range("AA27").Value = 10
$ws.Range("AB27").Value = 30
$ws.Range("AC27").Value = 20
$ws.Range("AD27").Value = 26
$ws.Range("AF27").Value = 7.1
$ws.Range("AG27").Value = 12.5
$ws.Range("AJ27").Value = 10
$ws.Range("AK27").Value = 13
$ws.Range("AL27").Value = 9.25
$ws.Range("AM27").Value = 24
$ws.Range("AN27").Value = 17
$ws.Range("AO27").Value = 23
$ws.Range("G27").Value = 2.7
$ws.Range("H27").Value = 3.6
$ws.Range("I27").Value = 2.32
$ws.Range("J27").Value = 3.2
$ws.Range("K27").Value = 2.25
$ws.Range("L27").Value = 2.85
$ws.Range("P27").Value = 3.85
$ws.Range("Y27").Value = 10.75
$ws.Range("Z27").Value = 15

# Row 28
$ws.Range("AA28").Value = 10.5
$ws.Range("AB28").Value = 35
$ws.Range("AC28").Value = 21
$ws.Range("AF28").Value = 7.4
$ws.Range("AJ28").Value = 10.25
$ws.Range("AK28").Value = 12.5
$ws.Range("AL28").Value = 9
$ws.Range("AM28").Value = 22
$ws.Range("AN28").Value = 15.5
$ws.Range("G28").Value = 2.9
$ws.Range("H28").Value = 3.65
$ws.Range("I28").Value = 2.18
$ws.Range("J28").Value = 3.3
$ws.Range("K28").Value = 2.3
$ws.Range("L28").Value = 2.67
$ws.Range("Q28").Value = 1.6
$ws.Range("R28").Value = 2.22
$ws.Range("X28").Value = 2.35
$ws.Range("Y28").Value = 12.5
$ws.Range("Z28").Value = 17.5
